$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53, shifting existing rows 53..92 down to 54..93
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new weekly price record
$ws.Cells.Item(53, 1).Value = 11
$ws.Cells.Item(53, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(53, 3).Value = "Bíobío"
$ws.Cells.Item(53, 4).Value = 44587
$ws.Cells.Item(53, 5).Value = 8
$ws.Cells.Item(53, 6).Value = 100112032
$ws.Cells.Item(53, 7).Value = "Zapallo italiano"
$ws.Cells.Item(53, 8).Value = "Sin especificar"
$ws.Cells.Item(53, 9).Value = "Primera"
$ws.Cells.Item(53, 10).Value = 310
$ws.Cells.Item(53, 11).Value = 10000
$ws.Cells.Item(53, 12).Value = 11000
$ws.Cells.Item(53, 13).Value = 10516
$ws.Cells.Item(53, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(53, 15).Value = "Región Metropolitana"
$ws.Cells.Item(53, 16).Value = 175
$ws.Cells.Item(53, 17).Value = 60
$ws.Cells.Item(53, 18).Value = "Hortaliza"
